$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97-118 down to 98-119.
$ws.Range("A97:R97").Insert()

# Copy formatting from the row above (row 96, now still row 96) down into the new row 97
# so that styles (e.g. date format on column D) are preserved, mirroring what Excel does
# automatically when inserting a row by copying the row above's format.
$ws.Range("A96:R96").Copy()
$ws.Range("A97:R97").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 97 values
$ws.Cells.Item(97, 1).Value = 1
$ws.Cells.Item(97, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(97, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(97, 4).Value = 44798
$ws.Cells.Item(97, 5).Value = 15
$ws.Cells.Item(97, 6).Value = 100112008
$ws.Cells.Item(97, 7).Value = "Coliflor"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Tercera"
$ws.Cells.Item(97, 10).Value = 1000
$ws.Cells.Item(97, 11).Value = 800
$ws.Cells.Item(97, 12).Value = 900
$ws.Cells.Item(97, 13).Value = 850
$ws.Cells.Item(97, 14).Value = "`$/unidad"
$ws.Cells.Item(97, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97, 16).Value = 850
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
